{"js": "// The document contains five \"<id>...</id>\" markers whose visible text was\n// split across three separate runs (the literal \"<id>\"/\"</id>\" tag text is\n// Courier New / color 7f6000 / size 9pt, while the identifier in between is\n// Arial / color 000000 / size 11pt), e.g.:\n//   \"<id>\" + \"p104r_a1\" + \"</id>\"\n// Each needs to become a single run reading \"<id>p104r_N</id>\" (the \"a\" is\n// dropped from \"p104r_aN\" -> \"p104r_N\"). Searching for the full visible\n// string and replacing it (Range.insertText(..., \"Replace\")) merges the\n// three runs into a single run that carries the formatting of the first\n// matched run, matching the target edit.\n\nconst body = context.document.body;\n\nfor (let i = 1; i <= 5; i++) {\n  const searchText = `<id>p104r_a${i}</id>`;\n  const replacementText = `<id>p104r_${i}</id>`;\n\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacementText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains five \"<id>...</id>\" markers whose visible text was\n# split across three separate runs (different fonts/colors for the literal\n# \"<id>\"/\"</id>\" tag text vs. the identifier itself), e.g.:\n#   \"<id>\" + \"p104r_a1\" + \"</id>\"\n# Each needs to become a single run reading \"<id>p104r_N</id>\" (the \"a\" is\n# dropped from \"p104r_aN\" -> \"p104r_N\"). Using Find & Replace on the full\n# visible string (\"<id>p104r_aN</id>\" -> \"<id>p104r_N</id>\") merges the three\n# runs into one run that carries the formatting of the first matched run,\n# which is exactly the observed result.\n\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le 5; $i++) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = \"<id>p104r_a$i</id>\"\n    $find.Replacement.Text = \"<id>p104r_$i</id>\"\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
